$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 4de3000e row (row 4, column G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-23 02:44:59"

# zh-cn sheet: "Correspond Handoff Datetime" (H4) and "Correspond Handback DateTime" (K4)
# for the 4de3000e row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-23 02:44:54"
$wsZhCn.Range("K4").Value = "2016-08-23 02:45:20"

# de-de sheet: "Correspond Handback DateTime" (K4) for the 4de3000e row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-23 02:45:27"
